# Changes of 5th April
# Refresh the FedEx shipment tracking report: new ShipmentTracking numbers
# (column P) for rows 2-25, with a handful of rows also getting updated
# ExpectedRate (column Q) and/or Result (column R) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric (tracking numbers, "$nn.nn"
# rate strings) while keeping the cell's native General/no-style format
# and text cell-type -- mirrors the source file where these columns are
# plain text, not numbers. We flip the format to Text, assign, then put
# the original style back so no lasting style/number-format change is
# left on the cell.
function Set-TextValue($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$data = @(
    @{Row=2; P="320018248263"; Q="`$19.04"; R="PASS"},
    @{Row=3; P="320018248274"; Q="`$27.50"; R="PASS"},
    @{Row=4; P="320018248300"; Q="`$31.73"; R="PASS"},
    @{Row=5; P="320018248322"; Q="`$43.36"; R="PASS"},
    @{Row=6; P="320018248366"; Q="`$56.05"; R="PASS"},
    @{Row=7; P="320018248388"; Q="`$231.08"; R="FAIL"},
    @{Row=8; P="320018248414"; Q="`$19.04"; R="PASS"},
    @{Row=9; P="320018248436"; Q="`$23.27"; R="PASS"},
    @{Row=10; P="320018248469"; Q="`$27.50"; R="PASS"},
    @{Row=11; P="320018248480"; Q="`$40.19"; R="PASS"},
    @{Row=12; P="320018248528"; Q="`$52.88"; R="PASS"},
    @{Row=13; P="320018248540"; Q="`$14.81"; R="PASS"},
    @{Row=14; P="320018248572"; Q="`$17.98"; R="PASS"},
    @{Row=15; P="320018248594"; Q="`$21.15"; R="PASS"},
    @{Row=16; P="320018248620"; Q="`$31.73"; R="PASS"},
    @{Row=17; P="320018248642"; Q="`$42.30"; R="PASS"},
    @{Row=18; P="320018248686"; Q="`$43.36"; R="PASS"},
    @{Row=19; P="320018248723"; Q="`$53.93"; R="PASS"},
    @{Row=20; P="320018248756"; Q="`$62.39"; R="PASS"},
    @{Row=21; P="320018248778"; Q="`$111.04"; R="PASS"},
    @{Row=22; P="320018248804"; Q="`$223.37"; R="FAIL"},
    @{Row=23; P="320018248815"; Q="`$436.98"; R="FAIL"},
    @{Row=24; P="320018248826"; Q="`$248.51"; R="FAIL"},
    @{Row=25; P="320018248837"; Q="`$52.88"; R="PASS"}
)

foreach ($item in $data) {
    Set-TextValue $ws "P$($item.Row)" $item.P
    Set-TextValue $ws "Q$($item.Row)" $item.Q
    $ws.Range("R$($item.Row)").Value = $item.R
}
